$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate rows 7-11 (rows 2-6 repeated verbatim in the original)
$ws.Range("A7:C11").Delete(-4162)

# Update the Genre column (B) with real-ish looking values
$ws.Range("B2").Value = "Genre not available"
$ws.Range("B3").Value = "Action, Drama"
$ws.Range("B4").Value = "Genre not available"
$ws.Range("B5").Value = "Genre not available"
$ws.Range("B6").Value = "Sci-Fi"

# Fill in the previously-missing Interstellar plot summary
$ws.Range("C4").Value = 'In the mid-21st century, humanity faces extinction due to dust storms and widespread crop blights. Joseph Cooper, a widowed former NASA test pilot, works as a farmer and raises his children, Murph and Tom, alongside his father-in-law Donald. Living in a post-truth society, Cooper is reprimanded by Murph''s teachers for telling her that the Apollo missions were not fabricated. During a dust storm, the two discover that dust patterns in Murph''s room, which she first attributes to a ghost, result from a gravitational anomaly, and translate into geographic coordinates. These lead them to a secret NASA facility headed by Professor John Brand, who explains that, 48 years earlier, a wormhole appeared near Saturn, leading to a system in another galaxy with twelve potentially habitable planets located near a black hole named Gargantua. Volunteers of the Lazarus expedition had previously travelled through the wormhole to evaluate the planets, with Miller, Edmunds, and Mann reporting back desirable results.'

# Add the new header columns (D:H), copying A1's header style so they match Name/Genre/Plot
$ws.Range("A1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Director"
$ws.Range("E1").Value = "Cast"
$ws.Range("F1").Value = "Release Date"
$ws.Range("G1").Value = "Runtime"
$ws.Range("H1").Value = "Rating"

# Fill the new columns with placeholder "Not found" values for every movie row
$ws.Range("D2:H6").Value = "Not found"
